$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.402.40"
$ws.Range("E2").Value = "  -5.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.276.27"
$ws.Range("E3").Value = "  -5.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.48"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.82"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.267.62"
$ws.Range("E8").Value = "  -5.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.31"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.813.99"
$ws.Range("E13").Value = "  -5.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.120"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000167"
$ws.Range("E15").Value = "  -5.05%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.253.64"
$ws.Range("E16").Value = "  -6.26%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.33"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.417.90"
$ws.Range("E18").Value = "  -5.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.65"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.04"
$ws.Range("E21").Value = "  -8.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "353.15"
$ws.Range("E22").Value = "  -7.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.391.57"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.31"
$ws.Range("E26").Value = "  -6.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000108"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.43"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.09"
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.150"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.292.00"
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.63"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.81"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "157.78"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0752"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.91"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.34"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.736"
$ws.Range("E45").Value = "  -7.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.55"
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.55"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.67"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.866"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.41"
$ws.Range("E51").Value = "  +5.59%  "
